$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: shift the dch / dchperceptron / KSVD data blocks down by copying
# whole rows (A:N) in an order that never overwrites data before it has been
# read (process the highest source row first).
# ---------------------------------------------------------------------------
$moves = @(
    @{src=13; dst=16},
    @{src=12; dst=15},
    @{src=11; dst=14},
    @{src=9;  dst=12},
    @{src=8;  dst=11},
    @{src=6;  dst=9},
    @{src=5;  dst=8}
)

foreach ($m in $moves) {
    $srcRange = $ws.Range("A" + $m.src + ":N" + $m.src)
    $dstRange = $ws.Range("A" + $m.dst + ":N" + $m.dst)
    $srcRange.Copy($dstRange)
}

# ---------------------------------------------------------------------------
# Step 2: row 13 becomes the new blank separator row (values cleared, style
# preserved).
# ---------------------------------------------------------------------------
$ws.Range("A13:N13").ClearContents()

# ---------------------------------------------------------------------------
# Step 3: rows 5 and 6 become the new "dl" algorithm rows (run did not
# terminate -> same placeholder pattern as the KSVD row that is now row 14).
# ---------------------------------------------------------------------------
$ws.Range("H14:L14").Copy($ws.Range("H5:L5"))
$ws.Range("H14:L14").Copy($ws.Range("H6:L6"))

# The placeholder pattern leaves J and L blank, but Copy() does not clear a
# destination cell when the matching source cell is blank - clear them
# explicitly instead.
$ws.Range("J5").ClearContents()
$ws.Range("L5").ClearContents()
$ws.Range("J6").ClearContents()
$ws.Range("L6").ClearContents()

$ws.Range("B5").Value = "dl"
$ws.Range("B6").Value = "dl"

$ws.Range("F5:G5").ClearContents()
$ws.Range("F6:G6").ClearContents()

$ws.Range("M5").Value = 2
$ws.Range("M6").Value = 2

# ---------------------------------------------------------------------------
# Step 4: reproduce the final selection recorded in the saved workbook.
# ---------------------------------------------------------------------------
$ws.Range("M19").Select()
